$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.391302333333333
$ws.Range("H2").Value = 4.173907
$ws.Range("I2").Value = 0.03171126955348368
$ws.Range("J2").Value = 0.03171126955348368
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.391302333333333
$ws.Range("N2").Value = 4.173907
$ws.Range("O2").Value = 0.03171126955348368
$ws.Range("P2").Value = 0.03171126955348368
$ws.Range("Q2").Value = 1.935722182738777
$ws.Range("R2").Value = 17.421499644649
$ws.Range("S2").Value = 0.001005604616693701
$ws.Range("T2").Value = 0.001005604616693701
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.391302333333333
$ws.Range("H3").Value = 4.173907
$ws.Range("I3").Value = 0.03171126955348368
$ws.Range("J3").Value = 0.03171126955348368
$ws.Range("O3").Value = 0.9045117298527411
$ws.Range("P3").Value = 0.9045117298527411
$ws.Range("Q3").Value = 55.21328678028377
$ws.Range("R3").Value = 496.919581022554
$ws.Range("S3").Value = 0.02868321527964808
$ws.Range("T3").Value = 0.02868321527964808
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.391302333333333
$ws.Range("H4").Value = 4.173907
$ws.Range("I4").Value = 0.03171126955348368
$ws.Range("J4").Value = 0.03171126955348368
$ws.Range("O4").Value = 0.06377700059377522
$ws.Range("P4").Value = 0.06377700059377524
$ws.Range("Q4").Value = 3.893081435598111
$ws.Range("R4").Value = 35.037732920383
$ws.Range("S4").Value = 0.002022449657141895
$ws.Range("T4").Value = 0.002022449657141895
$ws.Range("I5").Value = 0.9045117298527411
$ws.Range("J5").Value = 0.9045117298527411
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.391302333333333
$ws.Range("N5").Value = 4.173907
$ws.Range("O5").Value = 0.03171126955348368
$ws.Range("P5").Value = 0.03171126955348368
$ws.Range("Q5").Value = 55.21328678028377
$ws.Range("R5").Value = 496.919581022554
$ws.Range("S5").Value = 0.02868321527964808
$ws.Range("T5").Value = 0.02868321527964808
$ws.Range("I6").Value = 0.9045117298527411
$ws.Range("J6").Value = 0.9045117298527411
$ws.Range("O6").Value = 0.9045117298527411
$ws.Range("P6").Value = 0.9045117298527411
$ws.Range("S6").Value = 0.818141469441198
$ws.Range("T6").Value = 0.818141469441198
$ws.Range("I7").Value = 0.9045117298527411
$ws.Range("J7").Value = 0.9045117298527411
$ws.Range("O7").Value = 0.06377700059377522
$ws.Range("P7").Value = 0.06377700059377524
$ws.Range("S7").Value = 0.05768704513189492
$ws.Range("T7").Value = 0.05768704513189494
$ws.Range("I8").Value = 0.06377700059377522
$ws.Range("J8").Value = 0.06377700059377524
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.391302333333333
$ws.Range("N8").Value = 4.173907
$ws.Range("O8").Value = 0.03171126955348368
$ws.Range("P8").Value = 0.03171126955348368
$ws.Range("Q8").Value = 3.893081435598111
$ws.Range("R8").Value = 35.037732920383
$ws.Range("S8").Value = 0.002022449657141895
$ws.Range("T8").Value = 0.002022449657141895
$ws.Range("I9").Value = 0.06377700059377522
$ws.Range("J9").Value = 0.06377700059377524
$ws.Range("O9").Value = 0.9045117298527411
$ws.Range("P9").Value = 0.9045117298527411
$ws.Range("S9").Value = 0.05768704513189492
$ws.Range("T9").Value = 0.05768704513189494
$ws.Range("I10").Value = 0.06377700059377522
$ws.Range("J10").Value = 0.06377700059377524
$ws.Range("O10").Value = 0.06377700059377522
$ws.Range("P10").Value = 0.06377700059377524
$ws.Range("S10").Value = 0.004067505804738405
$ws.Range("T10").Value = 0.004067505804738407
